# Generate Report for Handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
# - Each localized sheet (zh-cn / de-de) gains a "Latest Target File" (E) and
#   "Latest Handback File" (F) column pair for the two tracked source files,
#   mirroring the existing "Source File Name" / "Latest Handoff File" link.
# - "Latest Handback DateTime" (G) is stamped and "Handoff Reason" (H) flips
#   from "Ignored" to "Include" for those same two rows.

$wb = $excel.ActiveWorkbook

$HYPERLINK_COLOR = 15570276  # BGR-encoded 0x6495ED -> matches the workbook's custom HyperLink style

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = $true
    $rng.Font.Color = $HYPERLINK_COLOR
}

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: Status mirrors (B/C columns) for both languages ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("B2").Value = $newStatus
$zhcn.Range("B3").Value = $newStatus

# Row 2: 206044d7-c91c-40dc-b682-66a175ceeef6
$zhcn.Range("E2").Value = "206044d7-c91c-40dc-b682-66a175ceeef6.md"
Style-AsHyperlink $zhcn.Range("E2")
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/beacc8af0d380dda163a7cc9282aaedaac687098/e2e/206044d7-c91c-40dc-b682-66a175ceeef6.md", "", "", "206044d7-c91c-40dc-b682-66a175ceeef6.md") | Out-Null

$zhcn.Range("F2").Value = "206044d7-c91c-40dc-b682-66a175ceeef6.aa7d2c2e4168c7ffb22ce69a9ac66a96353c629f.zh-cn.xlf"
Style-AsHyperlink $zhcn.Range("F2")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02417d0f3f3d8a047caa8a6e40c4f83d533ac06d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/206044d7-c91c-40dc-b682-66a175ceeef6.aa7d2c2e4168c7ffb22ce69a9ac66a96353c629f.zh-cn.xlf", "", "", "206044d7-c91c-40dc-b682-66a175ceeef6.aa7d2c2e4168c7ffb22ce69a9ac66a96353c629f.zh-cn.xlf") | Out-Null

$zhcn.Range("G2").Value = "2016-03-10 05:20:36"
$zhcn.Range("H2").Value = "Include"

# Row 3: e98baa6f-23ac-4935-b645-e2ada8ab1723
$zhcn.Range("E3").Value = "e98baa6f-23ac-4935-b645-e2ada8ab1723.md"
Style-AsHyperlink $zhcn.Range("E3")
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/beacc8af0d380dda163a7cc9282aaedaac687098/e2e/e98baa6f-23ac-4935-b645-e2ada8ab1723.md", "", "", "e98baa6f-23ac-4935-b645-e2ada8ab1723.md") | Out-Null

$zhcn.Range("F3").Value = "e98baa6f-23ac-4935-b645-e2ada8ab1723.725561f366f0a834cb29abffd96a28449324be5e.zh-cn.xlf"
Style-AsHyperlink $zhcn.Range("F3")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02417d0f3f3d8a047caa8a6e40c4f83d533ac06d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e98baa6f-23ac-4935-b645-e2ada8ab1723.725561f366f0a834cb29abffd96a28449324be5e.zh-cn.xlf", "", "", "e98baa6f-23ac-4935-b645-e2ada8ab1723.725561f366f0a834cb29abffd96a28449324be5e.zh-cn.xlf") | Out-Null

$zhcn.Range("G3").Value = "2016-03-10 05:20:36"
$zhcn.Range("H3").Value = "Include"

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("B2").Value = $newStatus
$dede.Range("B3").Value = $newStatus

# Row 2: 206044d7-c91c-40dc-b682-66a175ceeef6
$dede.Range("E2").Value = "206044d7-c91c-40dc-b682-66a175ceeef6.md"
Style-AsHyperlink $dede.Range("E2")
$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/beacc8af0d380dda163a7cc9282aaedaac687098/e2e/206044d7-c91c-40dc-b682-66a175ceeef6.md", "", "", "206044d7-c91c-40dc-b682-66a175ceeef6.md") | Out-Null

$dede.Range("F2").Value = "206044d7-c91c-40dc-b682-66a175ceeef6.aa7d2c2e4168c7ffb22ce69a9ac66a96353c629f.de-de.xlf"
Style-AsHyperlink $dede.Range("F2")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0d86c80d6da183b731ca0fb9147aa182189a663/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/206044d7-c91c-40dc-b682-66a175ceeef6.aa7d2c2e4168c7ffb22ce69a9ac66a96353c629f.de-de.xlf", "", "", "206044d7-c91c-40dc-b682-66a175ceeef6.aa7d2c2e4168c7ffb22ce69a9ac66a96353c629f.de-de.xlf") | Out-Null

$dede.Range("G2").Value = "2016-03-10 05:20:46"
$dede.Range("H2").Value = "Include"

# Row 3: e98baa6f-23ac-4935-b645-e2ada8ab1723
$dede.Range("E3").Value = "e98baa6f-23ac-4935-b645-e2ada8ab1723.md"
Style-AsHyperlink $dede.Range("E3")
$dede.Hyperlinks.Add($dede.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/beacc8af0d380dda163a7cc9282aaedaac687098/e2e/e98baa6f-23ac-4935-b645-e2ada8ab1723.md", "", "", "e98baa6f-23ac-4935-b645-e2ada8ab1723.md") | Out-Null

$dede.Range("F3").Value = "e98baa6f-23ac-4935-b645-e2ada8ab1723.725561f366f0a834cb29abffd96a28449324be5e.de-de.xlf"
Style-AsHyperlink $dede.Range("F3")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0d86c80d6da183b731ca0fb9147aa182189a663/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e98baa6f-23ac-4935-b645-e2ada8ab1723.725561f366f0a834cb29abffd96a28449324be5e.de-de.xlf", "", "", "e98baa6f-23ac-4935-b645-e2ada8ab1723.725561f366f0a834cb29abffd96a28449324be5e.de-de.xlf") | Out-Null

$dede.Range("G3").Value = "2016-03-10 05:20:46"
$dede.Range("H3").Value = "Include"

Write-Host "Handback report generated."
